# Update cryptos list with latest values (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.410.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.462.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  +1.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9515"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3656"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3064"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.79"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.035"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06584"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.420"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.142"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001025"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.462.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9705"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05840"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.433"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.428.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.075"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.616.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "112.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.848"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.898"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.07888"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7915"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.513"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05706"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.147"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.680"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02030"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9583"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.503"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1858"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5259"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.486"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "117.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5154"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.747"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06420"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9924"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.60%  "
